# Applies crypto price/volume updates and a few Coin/Link row swaps,
# as scraped on Mon Nov 11 05:49:45 UTC 2024 with GitHub Actions.
#
# Columns D (Price) and E (Volume(1h)) are stored as text in the sheet
# (e.g. "80.834.40", "  +1.55%  "), so force a text number format before
# assigning values to stop Excel auto-converting numeric-looking strings
# into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").NumberFormat = "@"
    $ws.Range("D2").Value = '80.834.40'
    $ws.Range("E2").NumberFormat = "@"
    $ws.Range("E2").Value = '  +1.55%  '
    $ws.Range("D3").NumberFormat = "@"
    $ws.Range("D3").Value = '3.134.68'
    $ws.Range("E3").NumberFormat = "@"
    $ws.Range("E3").Value = '  -2.45%  '
    $ws.Range("D4").NumberFormat = "@"
    $ws.Range("D4").Value = '0.999'
    $ws.Range("E4").NumberFormat = "@"
    $ws.Range("E4").Value = '  -0.04%  '
    $ws.Range("D5").NumberFormat = "@"
    $ws.Range("D5").Value = '205.53'
    $ws.Range("E5").NumberFormat = "@"
    $ws.Range("E5").Value = '  -1.25%  '
    $ws.Range("D6").NumberFormat = "@"
    $ws.Range("D6").Value = '623.23'
    $ws.Range("E6").NumberFormat = "@"
    $ws.Range("E6").Value = '  -1.51%  '
    $ws.Range("D7").NumberFormat = "@"
    $ws.Range("D7").Value = '0.282'
    $ws.Range("E7").NumberFormat = "@"
    $ws.Range("E7").Value = '  +23.41%  '
    $ws.Range("D8").NumberFormat = "@"
    $ws.Range("D8").Value = '0.999'
    $ws.Range("E8").NumberFormat = "@"
    $ws.Range("E8").Value = '  +0.01%  '
    $ws.Range("D9").NumberFormat = "@"
    $ws.Range("D9").Value = '0.576'
    $ws.Range("E9").NumberFormat = "@"
    $ws.Range("E9").Value = '  -1.97%  '
    $ws.Range("D10").NumberFormat = "@"
    $ws.Range("D10").Value = '3.130.76'
    $ws.Range("E10").NumberFormat = "@"
    $ws.Range("E10").Value = '  -2.59%  '
    $ws.Range("D11").NumberFormat = "@"
    $ws.Range("D11").Value = '0.572'
    $ws.Range("E11").NumberFormat = "@"
    $ws.Range("E11").Value = '  -2.10%  '
    $ws.Range("D12").NumberFormat = "@"
    $ws.Range("D12").Value = '0.0000250'
    $ws.Range("E12").NumberFormat = "@"
    $ws.Range("E12").Value = '  +10.13%  '
    $ws.Range("E13").NumberFormat = "@"
    $ws.Range("E13").Value = '  +0.74%  '
    $ws.Range("D14").NumberFormat = "@"
    $ws.Range("D14").Value = '5.25'
    $ws.Range("E14").NumberFormat = "@"
    $ws.Range("E14").Value = '  -3.43%  '
    $ws.Range("D15").NumberFormat = "@"
    $ws.Range("D15").Value = '3.713.05'
    $ws.Range("E15").NumberFormat = "@"
    $ws.Range("E15").Value = '  -2.27%  '
    $ws.Range("D16").NumberFormat = "@"
    $ws.Range("D16").Value = '31.11'
    $ws.Range("E16").NumberFormat = "@"
    $ws.Range("E16").Value = '  -2.48%  '
    $ws.Range("D17").NumberFormat = "@"
    $ws.Range("D17").Value = '80.833.50'
    $ws.Range("E17").NumberFormat = "@"
    $ws.Range("E17").Value = '  +1.95%  '
    $ws.Range("D18").NumberFormat = "@"
    $ws.Range("D18").Value = '3.145.49'
    $ws.Range("E18").NumberFormat = "@"
    $ws.Range("E18").Value = '  -1.49%  '
    $ws.Range("E19").NumberFormat = "@"
    $ws.Range("E19").Value = '  +9.99%  '
    $ws.Range("D20").NumberFormat = "@"
    $ws.Range("D20").Value = '13.87'
    $ws.Range("E20").NumberFormat = "@"
    $ws.Range("E20").Value = '  -4.37%  '
    $ws.Range("D21").NumberFormat = "@"
    $ws.Range("D21").Value = '429.94'
    $ws.Range("E21").NumberFormat = "@"
    $ws.Range("E21").Value = '  -0.85%  '
    $ws.Range("D22").NumberFormat = "@"
    $ws.Range("D22").Value = '8.91'
    $ws.Range("E22").NumberFormat = "@"
    $ws.Range("E22").Value = '  -5.95%  '
    $ws.Range("D23").NumberFormat = "@"
    $ws.Range("D23").Value = '5.05'
    $ws.Range("E23").NumberFormat = "@"
    $ws.Range("E23").Value = '  -0.44%  '
    $ws.Range("E24").NumberFormat = "@"
    $ws.Range("E24").Value = '  +4.20%  '
    $ws.Range("D25").NumberFormat = "@"
    $ws.Range("D25").Value = '5.13'
    $ws.Range("E25").NumberFormat = "@"
    $ws.Range("E25").Value = '  +7.33%  '
    $ws.Range("D26").NumberFormat = "@"
    $ws.Range("D26").Value = '3.308.13'
    $ws.Range("E26").NumberFormat = "@"
    $ws.Range("E26").Value = '  -1.69%  '
    $ws.Range("D27").NumberFormat = "@"
    $ws.Range("D27").Value = '75.46'
    $ws.Range("E27").NumberFormat = "@"
    $ws.Range("E27").Value = '  -2.48%  '
    $ws.Range("D28").NumberFormat = "@"
    $ws.Range("D28").Value = '10.83'
    $ws.Range("E28").NumberFormat = "@"
    $ws.Range("E28").Value = '  -1.74%  '
    $ws.Range("D29").NumberFormat = "@"
    $ws.Range("D29").Value = '0.997'
    $ws.Range("E29").NumberFormat = "@"
    $ws.Range("E29").Value = '  -0.80%  '
    $ws.Range("E30").NumberFormat = "@"
    $ws.Range("E30").Value = '  +4.54%  '
    $ws.Range("D31").NumberFormat = "@"
    $ws.Range("D31").Value = '0.999'
    $ws.Range("E31").NumberFormat = "@"
    $ws.Range("E31").Value = '  +0.01%  '
    $ws.Range("D32").NumberFormat = "@"
    $ws.Range("D32").Value = '8.94'
    $ws.Range("E32").NumberFormat = "@"
    $ws.Range("E32").Value = '  -0.35%  '
    $ws.Range("D33").NumberFormat = "@"
    $ws.Range("D33").Value = '553.56'
    $ws.Range("E33").NumberFormat = "@"
    $ws.Range("E33").Value = '  +6.47%  '
    $ws.Range("B34").Value = 'Fetch.AI'
    $ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    $ws.Range("D34").NumberFormat = "@"
    $ws.Range("D34").Value = '1.47'
    $ws.Range("E34").NumberFormat = "@"
    $ws.Range("E34").Value = '  -0.42%  '
    $ws.Range("B35").Value = 'Cronos'
    $ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    $ws.Range("D35").NumberFormat = "@"
    $ws.Range("D35").Value = '0.147'
    $ws.Range("E35").NumberFormat = "@"
    $ws.Range("E35").Value = '  +16.09%  '
    $ws.Range("B36").Value = 'PancakeSwap'
    $ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    $ws.Range("D36").NumberFormat = "@"
    $ws.Range("D36").Value = '1.99'
    $ws.Range("E36").NumberFormat = "@"
    $ws.Range("E36").Value = '  +0.22%  '
    $ws.Range("B37").Value = 'Kaspa'
    $ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    $ws.Range("D37").NumberFormat = "@"
    $ws.Range("D37").Value = '0.149'
    $ws.Range("E37").NumberFormat = "@"
    $ws.Range("E37").Value = '  +8.56%  '
    $ws.Range("D38").NumberFormat = "@"
    $ws.Range("D38").Value = '22.52'
    $ws.Range("E38").NumberFormat = "@"
    $ws.Range("E38").Value = '  -1.95%  '
    $ws.Range("D39").NumberFormat = "@"
    $ws.Range("D39").Value = '0.998'
    $ws.Range("E39").NumberFormat = "@"
    $ws.Range("E39").Value = '  -0.18%  '
    $ws.Range("D40").NumberFormat = "@"
    $ws.Range("D40").Value = '0.403'
    $ws.Range("E40").NumberFormat = "@"
    $ws.Range("E40").Value = '  -1.31%  '
    $ws.Range("E41").NumberFormat = "@"
    $ws.Range("E41").Value = '  +8.15%  '
    $ws.Range("E42").NumberFormat = "@"
    $ws.Range("E42").Value = '  +3.46%  '
    $ws.Range("D43").NumberFormat = "@"
    $ws.Range("D43").Value = '3.00'
    $ws.Range("E43").NumberFormat = "@"
    $ws.Range("E43").Value = '  +17.63%  '
    $ws.Range("B44").Value = 'Stacks'
    $ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    $ws.Range("D44").NumberFormat = "@"
    $ws.Range("D44").Value = '1.98'
    $ws.Range("E44").NumberFormat = "@"
    $ws.Range("E44").Value = '  +10.44%  '
    $ws.Range("B45").Value = 'Monero'
    $ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    $ws.Range("D45").NumberFormat = "@"
    $ws.Range("D45").Value = '160.23'
    $ws.Range("E45").NumberFormat = "@"
    $ws.Range("E45").Value = '  -2.33%  '
    $ws.Range("D47").NumberFormat = "@"
    $ws.Range("D47").Value = '185.79'
    $ws.Range("E47").NumberFormat = "@"
    $ws.Range("E47").Value = '  -5.91%  '
    $ws.Range("E48").NumberFormat = "@"
    $ws.Range("E48").Value = '  +0.72%  '
    $ws.Range("D49").NumberFormat = "@"
    $ws.Range("D49").Value = '43.69'
    $ws.Range("E49").NumberFormat = "@"
    $ws.Range("E49").Value = '  +1.21%  '
    $ws.Range("D50").NumberFormat = "@"
    $ws.Range("D50").Value = '0.768'
    $ws.Range("E50").NumberFormat = "@"
    $ws.Range("E50").Value = '  -5.09%  '
    $ws.Range("B51").Value = 'InjectiveProtocol'
    $ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    $ws.Range("D51").NumberFormat = "@"
    $ws.Range("D51").Value = '25.32'
    $ws.Range("E51").NumberFormat = "@"
    $ws.Range("E51").Value = '  +2.74%  '

